$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 11174
$ws1.Range("F5").Value = 1254
$ws1.Range("F6").Value = 1124
$ws1.Range("F8").Value = 296
$ws1.Range("F13").Value = 2153
$ws1.Range("F15").Value = 1060
$ws1.Range("F17").Value = 567
$ws1.Range("F19").Value = 962
$ws1.Range("F24").Value = 690
$ws1.Range("F28").Value = 55
$ws1.Range("F29").Value = 319
$ws1.Range("F31").Value = 187
$ws1.Range("F33").Value = 257
$ws1.Range("F34").Value = 604
$ws1.Range("F35").Value = 2184
$ws1.Range("F38").Value = 1481
$ws1.Range("F42").Value = 99

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 206
$ws2.Range("F7").Value = 76
$ws2.Range("F16").Value = 16

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 609

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1254
$ws4.Range("F7").Value = 609
$ws4.Range("F8").Value = 1124
$ws4.Range("F9").Value = 206
$ws4.Range("F10").Value = 296
$ws4.Range("F14").Value = 76
$ws4.Range("F16").Value = 2153
$ws4.Range("F18").Value = 1060
$ws4.Range("F20").Value = 567
$ws4.Range("F22").Value = 962
$ws4.Range("F27").Value = 690
$ws4.Range("F31").Value = 55
$ws4.Range("F33").Value = 187
$ws4.Range("F35").Value = 257
$ws4.Range("F36").Value = 2184
$ws4.Range("F40").Value = 1481
$ws4.Range("F44").Value = 99

$wb.Save()
